# Refresh the cryptocurrency price/volume snapshot (cryptos.xlsx) -
# GitHub Actions scheduled update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are preserved as exact text (matches original inline-string formatting)
$ws.Range("D2:D16").NumberFormat = "@"
$ws.Range("D18:D22").NumberFormat = "@"
$ws.Range("D24:D38").NumberFormat = "@"
$ws.Range("D40:D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns with the latest crypto snapshot
$ws.Range("D2").Value = "27.680.32"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "1.859.59"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("D4").Value = "1.037"
$ws.Range("E4").Value = "  +2.91%  "
$ws.Range("D5").Value = "323.16"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "1.032"
$ws.Range("E6").Value = "  +2.51%  "
$ws.Range("D7").Value = "0.4408"
$ws.Range("E7").Value = "  +2.72%  "
$ws.Range("D8").Value = "0.3800"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "0.07440"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").Value = "0.8815"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "21.71"
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("D12").Value = "1.868.41"
$ws.Range("E12").Value = "  -8.34%  "
$ws.Range("D13").Value = "5.550"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").Value = "6.734"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "0.07200"
$ws.Range("E15").Value = "  +4.12%  "
$ws.Range("D16").Value = "83.56"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "0.000009077"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").Value = "1.033"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").Value = "15.51"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").Value = "27.737.80"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("D22").Value = "5.303"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("E23").Value = "  +4.19%  "
$ws.Range("D24").Value = "158.55"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").Value = "1.930"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "18.82"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("D27").Value = "1.991"
$ws.Range("E27").Value = "  +4.71%  "
$ws.Range("D28").Value = "5.310"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "117.53"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").Value = "0.09082"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "1.209"
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("D32").Value = "0.7653"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").Value = "4.568"
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("D34").Value = "2.881"
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("D35").Value = "1.033"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("D36").Value = "1.156"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("D37").Value = "0.01983"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").Value = "0.05335"
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").Value = "2.830"
$ws.Range("E40").Value = "  +3.53%  "
$ws.Range("D41").Value = "0.1686"
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("D42").Value = "6.826"
$ws.Range("E42").Value = "  +5.86%  "
$ws.Range("D43").Value = "8.661"
$ws.Range("E43").Value = "  +4.63%  "
$ws.Range("D44").Value = "109.39"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("D45").Value = "10.58"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("D46").Value = "1.720"
$ws.Range("E46").Value = "  +4.10%  "
$ws.Range("D47").Value = "0.4676"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").Value = "0.06414"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").Value = "1.855"
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("D50").Value = "39.62"
$ws.Range("E50").Value = "  +4.52%  "
$ws.Range("D51").Value = "64.32"
$ws.Range("E51").Value = "  +1.18%  "
